$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before the "总计" (totals) sheet
# ---------------------------------------------------------------------
# An existing per-quarter fund-holding sheet already carries the bold /
# bordered / centered header style plus the index-column style we need to
# reuse for the new sheet (both are cellXfs index 2 in the original file).
$styleSource = $wb.Worksheets.Item("2021-Q4")

$total = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"

# NOTE: this runtime's worksheet references are positional anchors, so the
# old $total handle now resolves to the freshly-inserted sheet (same index
# slot). Re-resolve it by name after the structural change.
$total = $wb.Worksheets.Item("总计")

# Pull in the header-row / index-column formatting before writing any values.
$styleSource.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$styleSource.Range("A2").Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)

$q1headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $q1headers.Length; $c++) {
    $q1.Cells.Item(1, $c + 2).Value = $q1headers[$c]
}

$q1rows = @(
    @("513050", "易方达中证海外中国互联网50 QDII-ETF", "350.10", "98.05", "2.55", "8.9276", 8),
    @("164906", "交银施罗德中证海外中国互联网指数QDII-LOF", "110.11", "92.92", "3.54", "3.8979", 10),
    @("159605", "广发中证海外中国互联网30（QDII-ETF）", "29.04", "98.61", "5.43", "1.5769", 8),
    @("159607", "嘉实中证海外中国互联网30ETF（QDII）", "5.79", "98.25", "5.46", "0.3161", 8)
)

# The fund code / scale / position figures all look numeric, so they must be
# pre-formatted as Text before assignment - otherwise Excel's usual
# numeric-string coercion would turn them into real numbers.
$q1.Range("B2:B5").NumberFormat = "@"
$q1.Range("D2:G5").NumberFormat = "@"

for ($i = 0; $i -lt $q1rows.Length; $i++) {
    $r = $i + 2
    $q1.Cells.Item($r, 1).Value = $i
    $q1.Cells.Item($r, 2).Value = $q1rows[$i][0]
    $q1.Cells.Item($r, 3).Value = $q1rows[$i][1]
    $q1.Cells.Item($r, 4).Value = $q1rows[$i][2]
    $q1.Cells.Item($r, 5).Value = $q1rows[$i][3]
    $q1.Cells.Item($r, 6).Value = $q1rows[$i][4]
    $q1.Cells.Item($r, 7).Value = $q1rows[$i][5]
    $q1.Cells.Item($r, 8).Value = $q1rows[$i][6]
}

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" summary row into the "总计" sheet
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Restore the column-A index style (it gets stripped by ClearFormats/Insert)
$total.Range("A6").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 14.72

# Renumber the index column for the rows that shifted down
for ($r = 3; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}
